$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "id" values in column A for rows 2, 3 and 6 (leave other
# cells in those rows untouched - this is a ClearContents on individual
# cells, not a row/column delete, so nothing shifts).
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A6").ClearContents()

# Update the remembered selection on the sheet view.
$ws.Range("E15").Select() | Out-Null
